$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.477.41"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "1.803.98"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'224.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "'0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.30%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'39.62"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.77%  "
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("E10").Value = "  -3.48%  "
$ws.Range("D11").Value = "'0.0986"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Value = "2.064.90"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'11.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.85%  "
$ws.Range("D14").Value = "1.799.84"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "34.451.62"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'4.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").Value = "'239.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E20").Value = "  -3.12%  "
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "'4.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").Value = "'171.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "'17.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.95%  "
$ws.Range("D27").Value = "'7.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.97%  "
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("D33").Value = "'3.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -4.77%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "'1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.303.07"
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").Value = "'2.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.63%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "'81.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("D45").Value = "'13.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "1.965.93"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -5.07%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "'102.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0613"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
